$d = $word.ActiveDocument

# --- 1) Remove the "Tao là phúc" run, and merge away the now-empty
#        paragraph that followed it, leaving the (still bookmarked)
#        paragraph empty. ---
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.TrimEnd([char]13, [char]7) -eq "Tao là phúc") {
        # Delete just the run text (keep the paragraph mark + bookmark).
        $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $textRange.Delete()

        # Delete the following (empty) paragraph entirely, which merges
        # it into this one.
        $nextPara = $paras.Item($i + 1)
        $nextPara.Range.Delete()
        break
    }
}

# --- 2) Footer page-number field result: "1" -> "6". ---
$footer = $d.Sections.Item(1).Footers.Item(1)
$char1 = $footer.Range.Characters.Item(1)
$char1.Text = "6"
